# Refined metadata to be additional tab
#
# This script:
#   1. Updates the "panel_query_time" (time_taken) values on the "data"
#      sheet to reflect the new query run timestamps.
#   2. Adds a new "metadata" worksheet after "data" summarising the panel
#      query (data_name, data_id, data_version, data_version_created,
#      panel_query_time, panel_get_request).

$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Refresh the per-row query timestamps on the "data" sheet (column F)
# ---------------------------------------------------------------------
$newTimestamps = @{
    2 = "2021-10-05 14:21:18.236439"
    3 = "2021-10-05 14:21:18.236445"
    4 = "2021-10-05 14:21:18.236448"
    5 = "2021-10-05 14:21:18.236450"
    6 = "2021-10-05 14:21:18.236452"
    7 = "2021-10-05 14:21:18.236454"
    8 = "2021-10-05 14:21:18.236456"
    9 = "2021-10-05 14:21:18.236458"
    10 = "2021-10-05 14:21:18.236461"
    11 = "2021-10-05 14:21:18.236463"
    12 = "2021-10-05 14:21:18.236465"
    13 = "2021-10-05 14:21:18.236467"
    14 = "2021-10-05 14:21:18.236469"
    15 = "2021-10-05 14:21:18.236471"
    16 = "2021-10-05 14:21:18.236473"
    17 = "2021-10-05 14:21:18.236475"
    18 = "2021-10-05 14:21:18.236477"
    19 = "2021-10-05 14:21:18.236479"
    20 = "2021-10-05 14:21:18.236481"
    21 = "2021-10-05 14:21:18.236483"
    22 = "2021-10-05 14:21:18.236486"
    23 = "2021-10-05 14:21:18.236487"
    24 = "2021-10-05 14:21:18.236489"
    25 = "2021-10-05 14:21:18.236492"
    26 = "2021-10-05 14:21:18.236494"
    27 = "2021-10-05 14:21:18.236496"
    28 = "2021-10-05 14:21:18.236498"
    29 = "2021-10-05 14:21:18.236500"
    30 = "2021-10-05 14:21:18.236502"
    31 = "2021-10-05 14:21:18.236504"
    32 = "2021-10-05 14:21:18.236506"
    33 = "2021-10-05 14:21:18.236508"
    34 = "2021-10-05 14:21:18.236510"
}

foreach ($row in $newTimestamps.Keys) {
    $dataWs.Range("F" + $row).Value = $newTimestamps[$row]
}

# ---------------------------------------------------------------------
# 2. Add the new "metadata" worksheet after "data"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$metaWs = $wb.Worksheets.Add($null, $lastSheet)
$metaWs.Name = "metadata"

# Header row (bold, centered, bordered - matching the "data" sheet header)
$headerRange = $metaWs.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$metaWs.Range("B1").Value = "data_name"
$metaWs.Range("C1").Value = "data_id"
$metaWs.Range("D1").Value = "data_version"
$metaWs.Range("E1").Value = "data_version_created"
$metaWs.Range("F1").Value = "panel_query_time"
$metaWs.Range("G1").Value = "panel_get_request"

# Index cell (bold, centered, bordered - matching the "data" sheet index column)
$indexCell = $metaWs.Range("A2")
$indexCell.Font.Bold = $true
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160
$indexCell.Borders.LineStyle = 1
$indexCell.Value = 0

$metaWs.Range("B2").Value = "Intracerebral calcification disorders"
$metaWs.Range("C2").Value = 315
$metaWs.Range("D2").NumberFormat = "@"
$metaWs.Range("D2").Value = "1.28"
$metaWs.Range("E2").Value = "2021-05-10T15:43:20.495748Z"
$metaWs.Range("F2").Value = "2021-10-05 14:21:18.233991"
$metaWs.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/315/?format=json"

Write-Output "metadata sheet added; timestamps refreshed"
